# Update cryptos price/volume data per commit "Updated cryptos list on Thu Nov 28 04:39:05 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'95.776.42"
$ws.Range("E2").Value = '  +3.24%  '
$ws.Range("D3").Value = "'3.606.15"
$ws.Range("E3").Value = '  +5.32%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = "'240.86"
$ws.Range("E5").Value = '  +3.96%  '
$ws.Range("D6").Value = "'655.77"
$ws.Range("E6").Value = '  +5.81%  '
$ws.Range("E7").Value = '  +7.08%  '
$ws.Range("D8").Value = "'0.412"
$ws.Range("E8").Value = '  +5.40%  '
$ws.Range("E9").Value = '  -0.20%  '
$ws.Range("D10").Value = "'1.02"
$ws.Range("E10").Value = '  +5.88%  '
$ws.Range("D11").Value = "'3.603.19"
$ws.Range("E11").Value = '  +5.22%  '
$ws.Range("D12").Value = "'43.55"
$ws.Range("E12").Value = '  +1.69%  '
$ws.Range("D13").Value = "'0.201"
$ws.Range("E13").Value = '  +1.69%  '
$ws.Range("D14").Value = "'6.35"
$ws.Range("E14").Value = '  +1.79%  '
$ws.Range("D15").Value = "'4.255.10"
$ws.Range("E15").Value = '  +4.63%  '
$ws.Range("D16").Value = "'95.594.40"
$ws.Range("E16").Value = '  +3.10%  '
$ws.Range("D17").Value = "'0.0000258"
$ws.Range("E17").Value = '  +4.98%  '
$ws.Range("D18").Value = "'3.597.77"
$ws.Range("E18").Value = '  +5.06%  '
$ws.Range("D19").Value = "'7.97"
$ws.Range("E19").Value = '  -1.45%  '
$ws.Range("D20").Value = "'12.53"
$ws.Range("E20").Value = '  +8.28%  '
$ws.Range("D21").Value = "'18.20"
$ws.Range("E21").Value = '  +2.39%  '
$ws.Range("D22").Value = "'3.52"
$ws.Range("E22").Value = '  +6.21%  '
$ws.Range("D23").Value = "'0.490"
$ws.Range("E23").Value = '  +12.12%  '
$ws.Range("D24").Value = "'512.95"
$ws.Range("E24").Value = '  +3.14%  '
$ws.Range("D25").Value = "'0.0000199"
$ws.Range("E25").Value = '  +8.48%  '
$ws.Range("D26").Value = "'6.69"
$ws.Range("E26").Value = '  +2.88%  '
$ws.Range("D27").Value = "'97.13"
$ws.Range("E27").Value = '  +6.62%  '
$ws.Range("D28").Value = "'12.90"
$ws.Range("E28").Value = '  +7.99%  '
$ws.Range("D29").Value = "'3.20"
$ws.Range("E29").Value = '  +18.05%  '
$ws.Range("D30").Value = "'11.39"
$ws.Range("E30").Value = '  +1.07%  '
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = '  -0.44%  '
$ws.Range("D32").Value = "'0.140"
$ws.Range("E32").Value = '  +3.92%  '
$ws.Range("D33").Value = "'1.01"
$ws.Range("E33").Value = '  +0.81%  '
$ws.Range("D34").Value = "'0.177"
$ws.Range("E34").Value = '  +3.24%  '
$ws.Range("D35").Value = "'32.11"
$ws.Range("E35").Value = '  +8.65%  '
$ws.Range("D36").Value = "'0.565"
$ws.Range("E36").Value = '  +4.56%  '
$ws.Range("D37").Value = "'8.26"
$ws.Range("E37").Value = '  +10.49%  '
$ws.Range("D38").Value = "'567.91"
$ws.Range("E38").Value = '  +2.05%  '
$ws.Range("D39").Value = "'1.47"
$ws.Range("E39").Value = '  +6.11%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").Value = "'0.934"
$ws.Range("E40").Value = '  +2.34%  '
$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("D42").Value = "'0.151"
$ws.Range("E42").Value = '  +0.90%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = "'5.78"
$ws.Range("E43").Value = '  +5.89%  '
$ws.Range("B44").Value = 'ImmutableX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D44").Value = "'1.74"
$ws.Range("E44").Value = '  +0.88%  '
$ws.Range("D45").Value = "'23.79"
$ws.Range("E45").Value = '  +0.49%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = "'34.39"
$ws.Range("E46").Value = '  +37.74%  '
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").Value = "'2.27"
$ws.Range("E47").Value = '  +8.14%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").Value = "'0.0421"
$ws.Range("E48").Value = '  +3.74%  '
$ws.Range("D49").Value = "'54.37"
$ws.Range("E49").Value = '  +2.41%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = "'8.20"
$ws.Range("E50").Value = '  +2.92%  '
$ws.Range("B51").Value = 'MantraDAO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D51").Value = "'3.45"
$ws.Range("E51").Value = '  -6.33%  '
